$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '75.168.45'
$ws.Range("E2").Value = '  +1.18%  '

# Row 3
$ws.Range("D3").Value = '2.865.74'
$ws.Range("E3").Value = '  +10.69%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.97%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '189.41'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.43%  '

# Row 7
$ws.Range("E7").Value = '  -0.09%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.566'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.21%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.195'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.53%  '

# Row 10
$ws.Range("D10").Value = '2.862.63'
$ws.Range("E10").Value = '  +10.24%  '

# Row 11
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.163'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.41%  '

# Row 12
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.375'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.40%  '

# Row 13
$ws.Range("E13").Value = '  +3.79%  '

# Row 14
$ws.Range("D14").Value = '3.387.78'
$ws.Range("E14").Value = '  +10.74%  '

# Row 15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '75.181.38'
$ws.Range("E15").Value = '  +1.47%  '

# Row 16
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.26%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000191'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.86%  '

# Row 18
$ws.Range("D18").Value = '2.858.69'
$ws.Range("E18").Value = '  +10.58%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +13.95%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.34%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.39%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.68%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.44%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.07%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.34%  '

# Row 26
$ws.Range("E26").Value = '  +0.13%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.43%  '

# Row 28
$ws.Range("D28").Value = '3.008.43'
$ws.Range("E28").Value = '  +10.91%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.19%  '

# Row 30
$ws.Range("E30").Value = '  +12.57%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.39%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '534.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.76%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.42'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.71%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.99'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.82%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.25%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.18%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '20.40'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.89%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '162.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.50%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.31'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.16%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '184.58'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +24.45%  '

# Row 42
$ws.Range("E42").Value = '  +0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.85%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.346'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.03%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.06%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.14%  '

# Row 47
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.92%  '

# Row 48
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.82%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0862'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.90%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.582'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.25%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.25%  '
